$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 (Fighters Resilience): add/adjust bonus columns ---
$ws.Range("E11").Value = 0.00055
$ws.Range("G11").Value = 0.00075
$ws.Range("H11").Value = 0.0005

# --- Row 12 (Quick Feet): bump type id and move-timeout bonus ---
$ws.Range("B12").Value = 7
$ws.Range("I12").Value = 0.001

# --- Row 14 (Celestial Prayer): healing bonus per level ---
$ws.Range("F14").Value = 0.0015

# --- Row 15 (Soldiers Strength): damage/ac bonus per level ---
$ws.Range("E15").Value = 0.0015
$ws.Range("G15").Value = 0.0015

# --- Row 18 (Natures Insight): add healing/ac bonus columns ---
$ws.Range("F18").Value = 0.00075
$ws.Range("G18").Value = 0.00015

# --- Row 20 (Alchemy): lower max level, raise skill bonus per level ---
$ws.Range("D20").Value = 200
$ws.Range("L20").Value = 0.005

# --- Row 23: rename "Lust for Gold" skill to "Kingmanship" ---
$ws.Range("A23").Value = "Kingmanship"

# --- New Row 24: add "Hells Anvil" skill for Blacksmith class ---
$ws.Range("A24").Value = "Hells Anvil"
$ws.Range("B24").Value = 0
$ws.Range("C24").Value = "This skill only applies to blacksmiths. The more you level this skill the more defence you will get, you can get up to +200% Attack and +300% Defence."
$ws.Range("D24").Value = 999
$ws.Range("E24").Value = 0.001
$ws.Range("G24").Value = 0.003
$ws.Range("K24").Value = 1
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = "Blacksmith"
$ws.Range("N24").Value = 0
